# Update "paises" (countries) COVID dashboard sheet:
#  - refresh the "last updated" timestamp
#  - refresh per-country statistics for the countries that changed
#  - the sheet is kept sorted by "Casos totales" (column B) descending,
#    so re-sorting after updating the values reproduces the row
#    reshuffling seen for Ghana (jumps above Suiza/Kirguistan) and
#    Benin (jumps above Sierra Leona).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the "datos actualizados" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 30 de Julio de 2020 a las 03:13"

# --- helper: locate the data row whose column A matches a country name ---
function Find-CountryRow($name) {
    $lastRow = 219
    for ($r = 4; $r -le $lastRow; $r++) {
        if ($ws.Cells.Item($r, 1).Value2 -eq $name) {
            return $r
        }
    }
    return -1
}

# --- 2. Update per-country rows with refreshed figures ---

$updates = @(
    @{ Name = "Estados Unidos";  B = 4567750; C = 66973; D = 2239727; E = 2174303; F = 0; G = 1429; H = 153720 },
    @{ Name = "Canada";          B = 115470;  C = 476;   D = 100465;  E = 6088;    F = 0; G = 5;    H = 8917   },
    @{ Name = "Ghana";           B = 35142;   C = 736;   D = 31286;   E = 3681;    F = 0; G = 7;    H = 175    },
    @{ Name = "Benin";           B = 1805;    C = 35;    D = 1036;    E = 733;     F = 0; G = 1;    H = 36     },
    @{ Name = "Comoras";         B = 378;     C = 24;    D = 330;     E = 41;      F = 0; G = 0;    H = 7      },
    @{ Name = "Polinesia Francesa"; D = 62;   E = 0 }
)

foreach ($u in $updates) {
    $row = Find-CountryRow $u.Name
    if ($row -gt 0) {
        if ($u.ContainsKey("B")) { $ws.Cells.Item($row, 2).Value = $u.B }
        if ($u.ContainsKey("C")) { $ws.Cells.Item($row, 3).Value = $u.C }
        if ($u.ContainsKey("D")) { $ws.Cells.Item($row, 4).Value = $u.D }
        if ($u.ContainsKey("E")) { $ws.Cells.Item($row, 5).Value = $u.E }
        if ($u.ContainsKey("F")) { $ws.Cells.Item($row, 6).Value = $u.F }
        if ($u.ContainsKey("G")) { $ws.Cells.Item($row, 7).Value = $u.G }
        if ($u.ContainsKey("H")) { $ws.Cells.Item($row, 8).Value = $u.H }
    }
}

# --- 3. Re-sort the data range by "Casos totales" (column B) descending,
#        which naturally reshuffles rows whose totals overtook neighbours ---
$dataRange = $ws.Range("A4:H219")
$sortKey = $ws.Range("B4:B219")
$dataRange.Sort($sortKey, 2)
